$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - LinearRegression (label unchanged)
$ws.Range("B2").Value = 8946733865695072
$ws.Range("C2").Value = 8946733865695065
$ws.Range("D2").Value = 8946733865695068

# Row 3 - RandomForestRegressor (label unchanged)
$ws.Range("B3").Value = 18549871614683.46
$ws.Range("C3").Value = 0.02680288239301544
$ws.Range("D3").Value = 1483623511851067

# Row 4 - label changes from GradientBoostingRegressor to DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.03196617862036957
$ws.Range("C4").Value = 0.0320466541469181
$ws.Range("D4").Value = 388304397165756.8

# Row 5 - label changes from AdaBoostRegressor to MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 32035604183668.41
$ws.Range("C5").Value = 42383725269985.56
$ws.Range("D5").Value = 472128229414579.5
